$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date
$ws.Name = "Through 2021-11-23"

# Update row 13 ("November (through 11-xx)") with the 2021-12-01 data pull
$ws.Range("A13").Value = "November (through 11-23)"
$ws.Range("C13").Value = 22
$ws.Range("D13").Value = 0.0435
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 50
$ws.Range("G13").Value = 0.1228
$ws.Range("I13").Value = 87
$ws.Range("J13").Value = 0.0225
$ws.Range("L13").Value = 40
$ws.Range("M13").Value = 0.1304
$ws.Range("N13").Value = 6
$ws.Range("O13").Value = 34
$ws.Range("P13").Value = 0.15
$ws.Range("Q13").Value = 8
$ws.Range("R13").Value = 154
$ws.Range("S13").Value = 0.0494
$ws.Range("U13").Value = 155
$ws.Range("V13").Value = 0.019

# Update row 14 ("Total") with the 2021-12-01 data pull
$ws.Range("C14").Value = 248
$ws.Range("D14").Value = 0.1174
$ws.Range("E14").Value = 59
$ws.Range("F14").Value = 484
$ws.Range("G14").Value = 0.1087
$ws.Range("I14").Value = 736
$ws.Range("J14").Value = 0.0788
$ws.Range("L14").Value = 589
$ws.Range("M14").Value = 0.1089
$ws.Range("N14").Value = 54
$ws.Range("O14").Value = 468
$ws.Range("P14").Value = 0.1034
$ws.Range("Q14").Value = 62
$ws.Range("R14").Value = 1157
$ws.Range("S14").Value = 0.0509
$ws.Range("U14").Value = 1506
$ws.Range("V14").Value = 0.0593

$wb.Save()
